# calorimetry : input and output consistency : done
#
# The fitted-species tables (constants_evaluated, enthalpies_calculated,
# input_enthalpies) previously carried rows for the fixed/input components
# (H, L, OH) in addition to the actual fitted complexes (HL, H2L, HOH,
# HOHD). Those input-component rows are removed so the sheets only list
# the species that are actually solved for, and the enthalpy signs /
# deviation for the remaining rows are corrected.

$wb = $excel.ActiveWorkbook

# --- constants_evaluated: drop the H, L, OH rows (rows 2-4) -----------
$wsConst = $wb.Worksheets.Item("constants_evaluated")
$wsConst.Rows.Item(2).Delete()
$wsConst.Rows.Item(2).Delete()
$wsConst.Rows.Item(2).Delete()

# --- enthalpies_calculated: drop the H, L, OH rows (rows 2-4) ---------
$wsEnth = $wb.Worksheets.Item("enthalpies_calculated")
$wsEnth.Rows.Item(2).Delete()
$wsEnth.Rows.Item(2).Delete()
$wsEnth.Rows.Item(2).Delete()

# Remaining rows are now: 2=HL, 3=H2L, 4=HOH, 5=HOHD.
# Fix the sign of the calculated enthalpy values for HL and H2L.
$wsEnth.Cells.Item(2, 2).Value = 13.452049754684
$wsEnth.Cells.Item(3, 2).Value = 15.4099741863607
# HOHD's deviation is not available any more - clear it.
$wsEnth.Cells.Item(5, 3).ClearContents()

# --- input_enthalpies: drop the H, L, OH rows (rows 2-4) ---------------
$wsInEnth = $wb.Worksheets.Item("input_enthalpies")
$wsInEnth.Rows.Item(2).Delete()
$wsInEnth.Rows.Item(2).Delete()
$wsInEnth.Rows.Item(2).Delete()
